$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20 (existing rows 20-42 shift down to 21-43).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44729
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112013
$ws.Cells.Item(20, 7).Value = "Alcachofa"
$ws.Cells.Item(20, 8).Value = "Argentina(o)"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 200
$ws.Cells.Item(20, 11).Value = 17000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 17500
$ws.Cells.Item(20, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(20, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 16).Value = 438
$ws.Cells.Item(20, 17).Value = 40
$ws.Cells.Item(20, 18).Value = "Hortaliza"
